$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weight")
$divs = $wb.Worksheets.Item("Dividends")
$tl = $wb.Worksheets.Item("Timeline")

# ---------------------------------------------------------------------------
# 1. Existing table (rows 1-5): restyle column K to bold right-aligned 0.0000
#    (K1 swaps from the plain bold style to the bordered "d." header style;
#     K2:K5 pick up right alignment on top of their existing bold 0.0000 fmt)
# ---------------------------------------------------------------------------
$divs.Range("G2").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null

$ws.Range("K2").HorizontalAlignment = -4152
$ws.Range("K3").HorizontalAlignment = -4152
$ws.Range("K4").HorizontalAlignment = -4152
$ws.Range("K5").HorizontalAlignment = -4152

# new blank row 6 cell under the Ratio column, matching the 0.0000 bold style
$ws.Range("G5").Copy() | Out-Null
$ws.Range("G6").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# 2. Firm 1 / Firm 2 table (rows 8-14): updated inputs + new layout
# ---------------------------------------------------------------------------
$ws.Range("B9").Value = 0.12
$ws.Range("C9").Value = 0.25
$ws.Range("B10").Value = 0.1
$ws.Range("C10").Value = 0.2

# blank helper cells E8:E10 next to the firm table
$ws.Range("E10").Font.Bold = $true

$tl.Range("B3").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null
$ws.Range("E9").NumberFormat = "0.00%"

$ws.Range("E8").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. New "Unlevered / Relevered Beta" mini block in columns G:H
# ---------------------------------------------------------------------------
$ws.Range("G8").Value = "Volitility "
$ws.Range("G8").HorizontalAlignment = -4152
$ws.Range("H8").Value = 0.13

$ws.Range("G9").Value = "Rate "
$ws.Range("G9").HorizontalAlignment = -4152
$ws.Range("H9").Value = 0.03

$ws.Range("G10").Value = "Expected Return "
$ws.Range("G10").HorizontalAlignment = -4152
$ws.Range("H10").Value = 7.0000000000000007E-2

$ws.Range("G12").Value = "Cash "
$ws.Range("G12").HorizontalAlignment = -4152
$ws.Range("H12").Value = 125000
$ws.Range("H12").NumberFormat = "#,##0"

$divs.Range("G2").Copy() | Out-Null
$ws.Range("G13").PasteSpecial(-4122) | Out-Null
$ws.Range("G13").Value = "Borrow "

$ws.Range("C8").Copy() | Out-Null
$ws.Range("H13").PasteSpecial(-4122) | Out-Null
$ws.Range("H13").Value = 50000
$ws.Range("H13").NumberFormat = "#,##0"

$ws.Range("G14").Value = "Ratio"
$ws.Range("G14").HorizontalAlignment = -4152
$ws.Range("H14").Formula = "=(H12+H13)/H12"

$ws.Range("G16").Value = "a."
$ws.Range("G16").HorizontalAlignment = -4152
$ws.Range("H16").Formula = "=H9+H14*(H10-H9)"
$ws.Range("H16").NumberFormat = "0.00%"

$ws.Range("G17").Value = "b."
$ws.Range("G17").HorizontalAlignment = -4152
$ws.Range("H17").Formula = "=H8*H14"
$ws.Range("H17").NumberFormat = "0.00%"

# ---------------------------------------------------------------------------
# 4. New "Rate / Expected Return / Volatility / Sharpe Ratio" block in K:L
# ---------------------------------------------------------------------------
$ws.Range("K8").Value = "Expected Return "
$ws.Range("K8").HorizontalAlignment = -4152
$ws.Range("L8").Value = 0.14000000000000001

$ws.Range("K9").Value = "Volatility "
$ws.Range("K9").HorizontalAlignment = -4152
$ws.Range("L9").Value = 0.2

$divs.Range("G2").Copy() | Out-Null
$ws.Range("K10").PasteSpecial(-4122) | Out-Null
$ws.Range("K10").Value = "Rate "
$ws.Range("L10").Value = 3.7999999999999999E-2

$divs.Range("A3").Copy() | Out-Null
$ws.Range("K11").PasteSpecial(-4122) | Out-Null
$ws.Range("K11").Value = "a. Sharpe Ratio "
$ws.Range("L11").Formula = "=(L8-L10)/L9"

# ---------------------------------------------------------------------------
# 5. Column widths / selection
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(11).EntireColumn.AutoFit() | Out-Null

$ws.Range("K11").Select() | Out-Null
